$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1. Insert a new column before column B ("Week_Start_Date")
#    This shifts ASIN, MyForecast, Amazon Mean/P70/P80/P90 Forecast,
#    Product Title, and is_holiday_week one column to the right (B->C ... I->J).
$ws.Columns.Item(2).Insert()

# 2. New column B header + force the column to store values as text so the
#    week-start dates remain literal "YYYY-MM-DD" strings (matching source data)
#    instead of being auto-converted into date serial numbers.
$ws.Range("B1").Value = "Week_Start_Date"
$ws.Range("B2:B17").NumberFormat = "@"

# 3. Per-row updates: corrected week labels (no leading zero), the new
#    Week_Start_Date text value, and corrected MyForecast/Amazon Mean values.
$ws.Range("A2").Value = "W1"
$ws.Range("B2").Value = "2025-01-05"
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 13

$ws.Range("A3").Value = "W2"
$ws.Range("B3").Value = "2025-01-12"
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 15

$ws.Range("A4").Value = "W3"
$ws.Range("B4").Value = "2025-01-19"
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 15

$ws.Range("A5").Value = "W4"
$ws.Range("B5").Value = "2025-01-26"
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 17

$ws.Range("A6").Value = "W5"
$ws.Range("B6").Value = "2025-02-02"
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 9
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 16

$ws.Range("A7").Value = "W6"
$ws.Range("B7").Value = "2025-02-09"
$ws.Range("D7").Value = 9
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = 11
$ws.Range("H7").Value = 16

$ws.Range("A8").Value = "W7"
$ws.Range("B8").Value = "2025-02-16"
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 9
$ws.Range("G8").Value = 12
$ws.Range("H8").Value = 17

$ws.Range("A9").Value = "W8"
$ws.Range("B9").Value = "2025-02-23"
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 12
$ws.Range("H9").Value = 17

$ws.Range("A10").Value = "W9"
$ws.Range("B10").Value = "2025-03-02"
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 11
$ws.Range("H10").Value = 15

$ws.Range("A11").Value = "W10"
$ws.Range("B11").Value = "2025-03-09"
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = 11
$ws.Range("H11").Value = 16

$ws.Range("A12").Value = "W11"
$ws.Range("B12").Value = "2025-03-16"
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 8
$ws.Range("G12").Value = 11
$ws.Range("H12").Value = 16

$ws.Range("A13").Value = "W12"
$ws.Range("B13").Value = "2025-03-23"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 9
$ws.Range("G13").Value = 12
$ws.Range("H13").Value = 17

$ws.Range("A14").Value = "W13"
$ws.Range("B14").Value = "2025-03-30"
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = 16

$ws.Range("A15").Value = "W14"
$ws.Range("B15").Value = "2025-04-06"
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = 6
$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = 15

$ws.Range("A16").Value = "W15"
$ws.Range("B16").Value = "2025-04-13"
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 6
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 15

$ws.Range("A17").Value = "W16"
$ws.Range("B17").Value = "2025-04-20"
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 14

# 4. is_holiday_week (now column J) becomes a proper boolean FALSE, not 0/n
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 10).Value = $false
}

# 5. Summary sheet: corrected 16-week forecast total
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value = "113"
